$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.185.90'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').Value = '3.420.43'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').Value = '  +0.40%  '
$ws.Range('D5').Value = "'413.41"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('D6').Value = "'128.59"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').Value = "'0.619"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.57%  '
$ws.Range('D9').Value = "'0.720"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.58%  '
$ws.Range('D10').Value = "'0.138"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.67%  '
$ws.Range('D11').Value = "'42.64"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.06%  '
$ws.Range('D12').Value = "'9.16"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('D13').Value = '3.967.90'
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').Value = "'0.140"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').Value = "'0.0000212"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.92%  '
$ws.Range('D16').Value = "'20.40"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.75%  '
$ws.Range('D17').Value = '3.405.25'
$ws.Range('E17').Value = '  -1.56%  '
$ws.Range('E18').Value = '  +3.22%  '
$ws.Range('E19').Value = '  -1.57%  '
$ws.Range('D20').Value = '62.199.69'
$ws.Range('E20').Value = '  +0.44%  '
$ws.Range('D21').Value = "'463.29"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.45%  '
$ws.Range('D22').Value = "'90.45"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.00%  '
$ws.Range('E23').Value = '  +2.84%  '
$ws.Range('D24').Value = "'13.36"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.51%  '
$ws.Range('D25').Value = "'10.36"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +17.98%  '
$ws.Range('E26').Value = '  +1.55%  '
$ws.Range('D27').Value = "'32.84"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.45%  '
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('D29').Value = "'7.65"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.14%  '
$ws.Range('D30').Value = "'11.85"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.04%  '
$ws.Range('E31').Value = '  -3.57%  '
$ws.Range('E32').Value = '  -0.98%  '
$ws.Range('E33').Value = '  -2.29%  '
$ws.Range('D34').Value = "'40.61"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.45%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').Value = "'57.90"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.62%  '
$ws.Range('D37').Value = "'0.0485"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.30%  '
$ws.Range('D38').Value = "'1.00"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.28%  '
$ws.Range('E39').Value = '  +4.24%  '
$ws.Range('D40').Value = "'0.324"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.08%  '
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('D42').Value = "'3.32"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.64%  '
$ws.Range('D43').Value = "'144.65"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.49%  '
$ws.Range('E44').Value = '  +9.76%  '
$ws.Range('E45').Value = '  +4.81%  '
$ws.Range('D46').Value = "'4.31"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.55%  '
$ws.Range('D47').Value = "'2.47"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +19.30%  '
$ws.Range('D48').Value = "'16.38"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.63%  '
$ws.Range('D49').Value = "'22.11"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.27%  '
$ws.Range('D50').Value = '0.0₃0503'
$ws.Range('E50').Value = '  +22.41%  '
$ws.Range('D51').Value = "'110.06"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.21%  '
